$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

# Title: "Status Update: 12/7/2020" -> "Project Wrap Up 12/7/2020"
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Project Wrap Up 12/7/2020"

# Body: update wording in first and third bullets, and drop the last two
# paragraphs (blank spacer + "Overall..." summary).
$bodyShape = $s.Shapes.Item(2)
$bullet1 = "We were able to resolve or improve on of the failures from the 2nd. The VGA's ability to interface with the FPGA to display the squares was verified"
$bullet2 = "Code was included that makes it so that the moles themselves pop down after a random amount of time, as well as allowing for the handling of more than one mole being up at the same time"
$bullet3 = "The issue with switch inputs not catching the mole has been improved. While the switch will still occasionally to catch the mole, it occurs much less frequently now."
$bodyShape.TextFrame.TextRange.Text = $bullet1 + "`r" + $bullet2 + "`r" + $bullet3
